$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NATMI Icam1-Itgal LR-pair output: expand from a single target cluster per
# sending cluster to the full 3x3 sending x target cluster grid (ECs, FAPs, sCs),
# and refresh all computed statistics accordingly (per Dr Hou advice).

# Row 2: ECs -> ECs  (Icam1 -> Itgal)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Icam1"
$ws.Range("C2").Value = "Itgal"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 99.11651100000002
$ws.Range("H2").Value = 297.3495330000001
$ws.Range("I2").Value = 0.799346251215574
$ws.Range("J2").Value = 0.7993462512155741
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 21.087087
$ws.Range("N2").Value = 63.261261
$ws.Range("O2").Value = 0.9808818221599021
$ws.Range("P2").Value = 0.9808818221599021
$ws.Range("Q2").Value = 2090.078490593457
$ws.Range("R2").Value = 18810.70641534112
$ws.Range("S2").Value = 0.7840642074290191
$ws.Range("T2").Value = 0.7840642074290192

# Row 3: ECs -> FAPs  (Icam1 -> Itgal)
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Icam1"
$ws.Range("C3").Value = "Itgal"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 99.11651100000002
$ws.Range("H3").Value = 297.3495330000001
$ws.Range("I3").Value = 0.799346251215574
$ws.Range("J3").Value = 0.7993462512155741
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3230143333333333
$ws.Range("N3").Value = 0.9690430000000001
$ws.Range("O3").Value = 0.01502525635066456
$ws.Range("P3").Value = 0.01502525635066456
$ws.Range("Q3").Value = 32.016053722991
$ws.Range("R3").Value = 288.1444835069191
$ws.Range("S3").Value = 0.01201038233745671
$ws.Range("T3").Value = 0.01201038233745671

# Row 4: ECs -> sCs  (Icam1 -> Itgal)
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Icam1"
$ws.Range("C4").Value = "Itgal"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 99.11651100000002
$ws.Range("H4").Value = 297.3495330000001
$ws.Range("I4").Value = 0.799346251215574
$ws.Range("J4").Value = 0.7993462512155741
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.08799
$ws.Range("N4").Value = 0.26397
$ws.Range("O4").Value = 0.00409292148943331
$ws.Range("P4").Value = 0.004092921489433309
$ws.Range("Q4").Value = 8.721261802890002
$ws.Range("R4").Value = 78.49135622601001
$ws.Range("S4").Value = 0.00327166144909818
$ws.Range("T4").Value = 0.00327166144909818

# Row 5: FAPs -> ECs  (Icam1 -> Itgal)
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Icam1"
$ws.Range("C5").Value = "Itgal"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 23.582852
$ws.Range("H5").Value = 70.74855599999999
$ws.Range("I5").Value = 0.1901889417714845
$ws.Range("J5").Value = 0.1901889417714845
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 21.087087
$ws.Range("N5").Value = 63.261261
$ws.Range("O5").Value = 0.9808818221599021
$ws.Range("P5").Value = 0.9808818221599021
$ws.Range("Q5").Value = 497.293651832124
$ws.Range("R5").Value = 4475.642866489115
$ws.Range("S5").Value = 0.1865528757594772
$ws.Range("T5").Value = 0.1865528757594772

# Row 6: FAPs -> FAPs  (Icam1 -> Itgal)
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Icam1"
$ws.Range("C6").Value = "Itgal"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 23.582852
$ws.Range("H6").Value = 70.74855599999999
$ws.Range("I6").Value = 0.1901889417714845
$ws.Range("J6").Value = 0.1901889417714845
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.3230143333333333
$ws.Range("N6").Value = 0.9690430000000001
$ws.Range("O6").Value = 0.01502525635066456
$ws.Range("P6").Value = 0.01502525635066456
$ws.Range("Q6").Value = 7.617599216878666
$ws.Range("R6").Value = 68.558392951908
$ws.Range("S6").Value = 0.002857637605178169
$ws.Range("T6").Value = 0.002857637605178169

# Row 7: FAPs -> sCs  (Icam1 -> Itgal)
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Icam1"
$ws.Range("C7").Value = "Itgal"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 23.582852
$ws.Range("H7").Value = 70.74855599999999
$ws.Range("I7").Value = 0.1901889417714845
$ws.Range("J7").Value = 0.1901889417714845
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.08799
$ws.Range("N7").Value = 0.26397
$ws.Range("O7").Value = 0.00409292148943331
$ws.Range("P7").Value = 0.004092921489433309
$ws.Range("Q7").Value = 2.07505514748
$ws.Range("R7").Value = 18.67549632732
$ws.Range("S7").Value = 0.0007784284068290893
$ws.Range("T7").Value = 0.0007784284068290891

# Row 8: sCs -> ECs  (Icam1 -> Itgal)
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Icam1"
$ws.Range("C8").Value = "Itgal"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.297604333333333
$ws.Range("H8").Value = 3.892813
$ws.Range("I8").Value = 0.01046480701294141
$ws.Range("J8").Value = 0.01046480701294141
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 21.087087
$ws.Range("N8").Value = 63.261261
$ws.Range("O8").Value = 0.9808818221599021
$ws.Range("P8").Value = 0.9808818221599021
$ws.Range("Q8").Value = 27.362695468577
$ws.Range("R8").Value = 246.264259217193
$ws.Range("S8").Value = 0.01026473897140569
$ws.Range("T8").Value = 0.01026473897140569

# Row 9: sCs -> FAPs  (Icam1 -> Itgal)
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Icam1"
$ws.Range("C9").Value = "Itgal"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.297604333333333
$ws.Range("H9").Value = 3.892813
$ws.Range("I9").Value = 0.01046480701294141
$ws.Range("J9").Value = 0.01046480701294141
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.3230143333333333
$ws.Range("N9").Value = 0.9690430000000001
$ws.Range("O9").Value = 0.01502525635066456
$ws.Range("P9").Value = 0.01502525635066456
$ws.Range("Q9").Value = 0.4191447986621111
$ws.Range("R9").Value = 3.772303187959
$ws.Range("S9").Value = 0.0001572364080296769
$ws.Range("T9").Value = 0.0001572364080296769

# Row 10: sCs -> sCs  (Icam1 -> Itgal)
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Icam1"
$ws.Range("C10").Value = "Itgal"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.297604333333333
$ws.Range("H10").Value = 3.892813
$ws.Range("I10").Value = 0.01046480701294141
$ws.Range("J10").Value = 0.01046480701294141
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.08799
$ws.Range("N10").Value = 0.26397
$ws.Range("O10").Value = 0.00409292148943331
$ws.Range("P10").Value = 0.004092921489433309
$ws.Range("Q10").Value = 0.11417620529
$ws.Range("R10").Value = 1.02758584761
$ws.Range("S10").Value = 0.00004283163350604029
$ws.Range("T10").Value = 0.00004283163350604028
